$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout (rows 1-10):
#  1: element | type
#  2: demand1 | demand
#  3: net1 | net
#  4: pv1 | pv
#  5: bat1 | bat
#  6: CHP1 | CHP              <- remove
#  7: solar_th1 | solar_th    <- remove
#  8: pvt1 | pvt
#  9: charging_station1 | charging_station
# 10: charging_station2 | charging_station  <- remove (becomes heat_pump1/heat_pump replacing row 9's duplicate)

# Remove the CHP1/CHP row (row 6) and solar_th1/solar_th row (originally row 7;
# after the first deletion it becomes row 6, so we delete it again).
$ws.Rows(6).Delete()
$ws.Rows(6).Delete()

# Now the layout is:
#  1: element | type
#  2: demand1 | demand
#  3: net1 | net
#  4: pv1 | pv
#  5: bat1 | bat
#  6: pvt1 | pvt
#  7: charging_station1 | charging_station
#  8: charging_station2 | charging_station

# Replace the last row's data with the new heat_pump entry instead of the
# duplicate charging_station2 row.
$ws.Range("A8").Value = "heat_pump1"
$ws.Range("B8").Value = "heat_pump"
